# Add a new "2022-Q3" sheet right after "总计", carrying the same layout
# as the other quarterly sheets, and insert a matching summary row on
# the "总计" sheet.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# --- 1. Create the new "2022-Q3" sheet right after "总计" -----------------
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Style = "Normal"

$q3.Range("A2").Value = 0
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "012346"
$q3.Range("C2").Value = "易方达港股通成长混合A"
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "20.23"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "85.52"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "2.89"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.5846"
$q3.Range("H2").Value = 9
$q3.Range("B2:H2").Style = "Normal"

$q3.Range("A3").Value = 1
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "012347"
$q3.Range("C3").Value = "易方达港股通成长混合C"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "6.17"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "85.52"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "2.89"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.1783"
$q3.Range("H3").Value = 9
$q3.Range("B3:H3").Style = "Normal"

# Headers (B1:H1) and the A-column index cells carry the bold/boxed style
# used throughout the workbook ("总计"!A2 is a convenient donor cell).
$totalSheet.Range("A2").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2:A3").PasteSpecial(-4122)

$q3.Range("A2").Value = 0
$q3.Range("A3").Value = 1

# --- 2. Insert the matching "2022-Q3" row into "总计" ----------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.76

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

# --- 3. Restore the original active sheet/selection -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
